# Updates currentAveragePrice / LevePrice / LeveProfit columns across Leve tables
# (scheduled-runner market-price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1338.6626
$ws.Range("I15").Value = 1338.6626
$ws.Range("K15").Value = 4015.9878
$ws.Range("M15").Value = -3846.9878
$ws.Range("H18").Value = 1104
$ws.Range("I18").Value = 814.5
$ws.Range("K18").Value = 814.5
$ws.Range("M18").Value = -530.5
$ws.Range("H76").Value = 7181.4443
$ws.Range("J76").Value = 7690.2
$ws.Range("L76").Value = 7690.2
$ws.Range("N76").Value = -8320.200000000001
$ws.Range("H79").Value = 7181.4443
$ws.Range("J79").Value = 7690.2
$ws.Range("L79").Value = 7690.2
$ws.Range("N79").Value = -9874.200000000001
$ws.Range("H87").Value = 41900
$ws.Range("I87").Value = 9833.333000000001
$ws.Range("K87").Value = 9833.333000000001
$ws.Range("M87").Value = -8585.333000000001
$ws.Range("H90").Value = 41900
$ws.Range("I90").Value = 9833.333000000001
$ws.Range("K90").Value = 29499.999
$ws.Range("M90").Value = -23259.999
$ws.Range("H92").Value = 3768.6667
$ws.Range("I92").Value = 2483.4443
$ws.Range("J92").Value = 5696.5
$ws.Range("K92").Value = 2483.4443
$ws.Range("L92").Value = 5696.5
$ws.Range("M92").Value = -1235.4443
$ws.Range("N92").Value = -8192.5
$ws.Range("H98").Value = 172552.94
$ws.Range("I98").Value = 798.42426
$ws.Range("K98").Value = 798.42426
$ws.Range("M98").Value = 699.57574
$ws.Range("H122").Value = 172552.94
$ws.Range("I122").Value = 798.42426
$ws.Range("K122").Value = 2395.27278
$ws.Range("M122").Value = 54.72721999999976
$ws.Range("H132").Value = 1746.7609
$ws.Range("I132").Value = 1751.4318
$ws.Range("J132").Value = 1644
$ws.Range("K132").Value = 5254.2954
$ws.Range("L132").Value = 4932
$ws.Range("M132").Value = -2724.2954
$ws.Range("N132").Value = -9992
$ws.Range("H137").Value = 12823626
$ws.Range("I137").Value = 47621770
$ws.Range("J137").Value = 3257
$ws.Range("K137").Value = 142865310
$ws.Range("L137").Value = 9771
$ws.Range("M137").Value = -142862760
$ws.Range("N137").Value = -14871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 283.33334
$ws.Range("I4").Value = 283.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 283.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -167.33334
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 392.58334
$ws.Range("J5").Value = 479.2
$ws.Range("L5").Value = 479.2
$ws.Range("N5").Value = -703.2
$ws.Range("H32").Value = 8312.311
$ws.Range("I32").Value = 7305.375
$ws.Range("K32").Value = 7305.375
$ws.Range("M32").Value = -7018.375
$ws.Range("H45").Value = 2639.5957
$ws.Range("I45").Value = 2373.5386
$ws.Range("J45").Value = 3936.625
$ws.Range("K45").Value = 2373.5386
$ws.Range("L45").Value = 3936.625
$ws.Range("M45").Value = -1996.5386
$ws.Range("N45").Value = -4690.625
$ws.Range("H61").Value = 4451.0977
$ws.Range("I61").Value = 3982.6
$ws.Range("K61").Value = 3982.6
$ws.Range("M61").Value = -3770.6
$ws.Range("H97").Value = 1755.3636
$ws.Range("I97").Value = 1430.9
$ws.Range("K97").Value = 1430.9
$ws.Range("M97").Value = -934.9000000000001
$ws.Range("H131").Value = 76799.664
$ws.Range("J131").Value = 76799.664
$ws.Range("L131").Value = 76799.664
$ws.Range("N131").Value = -86879.664
$ws.Range("H132").Value = 2668.962
$ws.Range("I132").Value = 2227.746
$ws.Range("K132").Value = 6683.238
$ws.Range("M132").Value = -4153.238
$ws.Range("H136").Value = 4451.0977
$ws.Range("I136").Value = 3982.6
$ws.Range("K136").Value = 11947.8
$ws.Range("M136").Value = -9397.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 392.58334
$ws.Range("J4").Value = 479.2
$ws.Range("L4").Value = 479.2
$ws.Range("N4").Value = -709.2
$ws.Range("H86").Value = 4925.0435
$ws.Range("I86").Value = 4067.7778
$ws.Range("J86").Value = 8011.2
$ws.Range("K86").Value = 4067.7778
$ws.Range("L86").Value = 8011.2
$ws.Range("M86").Value = -2944.7778
$ws.Range("N86").Value = -10257.2
$ws.Range("H89").Value = 4925.0435
$ws.Range("I89").Value = 4067.7778
$ws.Range("J89").Value = 8011.2
$ws.Range("K89").Value = 20338.889
$ws.Range("L89").Value = 40056
$ws.Range("M89").Value = -14722.889
$ws.Range("N89").Value = -51288
$ws.Range("H94").Value = 2648.6553
$ws.Range("I94").Value = 2642.44
$ws.Range("K94").Value = 2642.44
$ws.Range("M94").Value = -2191.44
$ws.Range("H99").Value = 1578.909
$ws.Range("I99").Value = 1333.4375
$ws.Range("J99").Value = 2233.5
$ws.Range("K99").Value = 1333.4375
$ws.Range("L99").Value = 2233.5
$ws.Range("M99").Value = 164.5625
$ws.Range("N99").Value = -5229.5
$ws.Range("H105").Value = 20317.945
$ws.Range("I105").Value = 20913.273
$ws.Range("K105").Value = 20913.273
$ws.Range("M105").Value = -19166.273
$ws.Range("H134").Value = 2733.389
$ws.Range("I134").Value = 2023.7587
$ws.Range("J134").Value = 5673.2856
$ws.Range("K134").Value = 6071.2761
$ws.Range("L134").Value = 17019.8568
$ws.Range("M134").Value = -3536.2761
$ws.Range("N134").Value = -22089.8568
$ws.Range("H140").Value = 66798.09
$ws.Range("J140").Value = 66798.09
$ws.Range("L140").Value = 66798.09
$ws.Range("N140").Value = -77158.09

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 434.125
$ws.Range("I7").Value = 199.14285
$ws.Range("J7").Value = 616.8889
$ws.Range("K7").Value = 199.14285
$ws.Range("L7").Value = 616.8889
$ws.Range("M7").Value = -86.14285000000001
$ws.Range("N7").Value = -842.8889
$ws.Range("H58").Value = 3598.6562
$ws.Range("I58").Value = 1380.0769
$ws.Range("K58").Value = 1380.0769
$ws.Range("M58").Value = -1177.0769
$ws.Range("H99").Value = 3143.0833
$ws.Range("J99").Value = 3336.6667
$ws.Range("L99").Value = 3336.6667
$ws.Range("N99").Value = -6332.6667
$ws.Range("H126").Value = 3143.0833
$ws.Range("J126").Value = 3336.6667
$ws.Range("L126").Value = 10010.0001
$ws.Range("N126").Value = -14950.0001
$ws.Range("H136").Value = 3598.6562
$ws.Range("I136").Value = 1380.0769
$ws.Range("K136").Value = 4140.2307
$ws.Range("M136").Value = -1590.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 64057.36
$ws.Range("J2").Value = 80059.45
$ws.Range("L2").Value = 480356.7
$ws.Range("N2").Value = -480582.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 452.95834
$ws.Range("J2").Value = 956.6
$ws.Range("L2").Value = 956.6
$ws.Range("N2").Value = -1182.6
$ws.Range("H58").Value = 25896.6
$ws.Range("I58").Value = 19994.666
$ws.Range("J58").Value = 34749.5
$ws.Range("K58").Value = 19994.666
$ws.Range("L58").Value = 34749.5
$ws.Range("M58").Value = -19717.666
$ws.Range("N58").Value = -35303.5
$ws.Range("H122").Value = 4609.6895
$ws.Range("I122").Value = 3264.9092
$ws.Range("K122").Value = 9794.7276
$ws.Range("M122").Value = -7344.7276
$ws.Range("H132").Value = 1645.6923
$ws.Range("I132").Value = 1201.9445
$ws.Range("K132").Value = 3605.8335
$ws.Range("M132").Value = -1075.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6038.3335
$ws.Range("I22").Value = 1499
$ws.Range("K22").Value = 1499
$ws.Range("M22").Value = -1204
$ws.Range("H25").Value = 3003.3333
$ws.Range("I25").Value = 505
$ws.Range("J25").Value = 8000
$ws.Range("K25").Value = 505
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = -275
$ws.Range("N25").Value = -8460
$ws.Range("H27").Value = 6038.3335
$ws.Range("I27").Value = 1499
$ws.Range("K27").Value = 1499
$ws.Range("M27").Value = -1392
$ws.Range("H61").Value = 6746.933
$ws.Range("I61").Value = 6110
$ws.Range("J61").Value = 8020.8
$ws.Range("K61").Value = 6110
$ws.Range("L61").Value = 8020.8
$ws.Range("M61").Value = -5908
$ws.Range("N61").Value = -8424.799999999999
$ws.Range("H113").Value = 6746.933
$ws.Range("I113").Value = 6110
$ws.Range("J113").Value = 8020.8
$ws.Range("K113").Value = 6110
$ws.Range("L113").Value = 8020.8
$ws.Range("M113").Value = -3940
$ws.Range("N113").Value = -12360.8
$ws.Range("H131").Value = 59994
$ws.Range("J131").Value = 59994
$ws.Range("L131").Value = 59994
$ws.Range("N131").Value = -70074

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()
